$wb = $excel.ActiveWorkbook

# Source of the existing bold/boxed header-ish style already used on
# column-A "index" cells elsewhere in the workbook (barnehage!A2 carries
# it). Copy/paste-special of this cell's format is used below so the new
# rows' A-column cells land on the very same shared style record instead
# of Excel/the engine minting a brand-new (duplicate) one.
$styleSrc = $wb.Worksheets.Item("barnehage").Range("A2")

# --- Sheet "foresatt" ---
$ws = $wb.Worksheets.Item("foresatt")

$styleSrc.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 2
$ws.Cells.Item(2, 3).Value = "erfan sarwari"
$ws.Cells.Item(2, 4).Value = "eigemyrveien 27"
$ws.Cells.Item(2, 5).Value = 93097239
$ws.Cells.Item(2, 6).Value = 30090578123

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = "Erfan sarwari"
$ws.Cells.Item(3, 4).Value = "eigemyrveien 26"
$ws.Cells.Item(3, 5).Value = 94086394
$ws.Cells.Item(3, 6).Value = 30090714256

# --- Sheet "barn" ---
$ws = $wb.Worksheets.Item("barn")

$styleSrc.Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 30090234152

# --- Sheet "soknad" ---
$ws = $wb.Worksheets.Item("soknad")

$styleSrc.Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = "on"
$ws.Cells.Item(2, 7).Value = "on"
$ws.Cells.Item(2, 8).Value = "on"

# fr_annet (I2) is present in the source row but blank; force the cell to
# exist without leaving any value/type behind.
$ws.Cells.Item(2, 9).NumberFormat = "@"
$ws.Cells.Item(2, 9).ClearFormats()

$ws.Cells.Item(2, 10).Value = "1, 2, 3, 4, 5"

# sosken__i_barnehagen (K2) likewise present but blank.
$ws.Cells.Item(2, 11).NumberFormat = "@"
$ws.Cells.Item(2, 11).ClearFormats()

# tidspunkt_oppstart (L2) must stay the literal text "2024-11-29" rather
# than be auto-converted into a date serial number, so force the cell to
# Text before typing it in, then drop the number-format again.
$ws.Range("L2").NumberFormat = "@"
$ws.Cells.Item(2, 12).Value = "2024-11-29"
$ws.Range("L2").ClearFormats()

$ws.Cells.Item(2, 13).Value = 900000
